$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D and E (ownTeam, oppTeam), shifting existing D:I to F:K
$ws.Range("D:E").EntireColumn.Insert()

# Set header values
$ws.Range("D1").Value = "ownTeam"
$ws.Range("E1").Value = "oppTeam"

# Set data values
$ws.Range("D2").Value = "Kings XI Punjab"
$ws.Range("E2").Value = "Sunrisers Hyderabad"

$wb.Save()
